$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LEAVE CREDITS")
$ws.Rows.Item(20).Insert()
$ws.Range("A10:K10").Copy()
$ws.Range("A20:K20").PasteSpecial(-4122)
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "2024"
$ws.Range("A20").NumberFormat = "mm/dd/yy;@"
"done"
